$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing score for row 2 (test tean 1)
$ws.Range("C2").Value = 3.6

# Add a new row for team "团队1"
$ws.Range("A3").Value = 33
$ws.Range("B3").Value = "团队1"
$ws.Range("C3").Value = 9.9
